$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells (I1 "I0", J1 "IF"), copying the existing
# header formatting (bold/border/centered) from H1 so the new columns match
# the look of the rest of the header row.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new I (I0) and J (IF) columns for the data rows (2-30).
# For rows 2-29, I is constant (1) and J mirrors the existing column H
# value for that row. Row 30 carries its own literal values.
$iValues = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,4)
$jValues = @(4,6,5,7,5,5,5,5,4,3,5,5,5,9,5,6,5,6,7,5,3,5,6,8,6,6,5,3,6)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
